# "novas imagens e resultados"
#
# Adds the new result percentages (columns M/N, rows 14-27) to Planilha1,
# and updates the view/selection state so that Planilha1 (not Planilha2)
# is the active tab, with R31 selected on Planilha1 and M5 (a single
# cell, not the old M5:R5 block) selected on Planilha2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# ---------------------------------------------------------------------
# New result percentages for rows 14-27, columns M and N (0.00% format)
# ---------------------------------------------------------------------
$newResults = @{
  14 = @(0.2599, 0.2454)
  15 = @(0.3277, 0.3471)
  16 = @(0.2629, 0.3105)
  17 = @(0.2024, 0.2648)
  18 = @(0.1101, 0.1151)
  19 = @(0.2299, 0.3990)
  20 = @(0.2193, 0.2268)
  21 = @(0.4112, 0.4714)
  22 = @(0.1609, 0.1815)
  23 = @(0.3184, 0.3782)
  24 = @(0.3319, 0.3454)
  25 = @(0.5561, 0.6226)
  26 = @(0.2507, 0.2081)
  27 = @(0.3705, 0.3784)
}

foreach ($row in $newResults.Keys) {
  $values = $newResults[$row]

  $mCell = $ws1.Range("M$row")
  $mCell.Value = $values[0]
  $mCell.NumberFormat = "0.00%"

  $nCell = $ws1.Range("N$row")
  $nCell.Value = $values[1]
  $nCell.NumberFormat = "0.00%"
}

# ---------------------------------------------------------------------
# View / selection state:
#  - Planilha2 becomes inactive and its selection collapses to the
#    single cell M5 (was M5:R5)
#  - Planilha1 becomes the active/selected tab with R31 selected
# ---------------------------------------------------------------------
[void]($ws2.Activate())
[void]($excel.ActiveWindow.DisplayGridlines = $true)
[void]($ws2.Range("M5").Select())

[void]($ws1.Activate())
[void]($excel.ActiveWindow.DisplayGridlines = $true)
[void]($ws1.Range("R31").Select())
